# "Generate Report for Handback"
# Refreshes the handoff/handback/generate timestamps for the
# a1611c47-fe89-4e11-a2a2-904b596ac550 row (row 3) across the Overview,
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-09-05 06:55:24"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-09-05 06:55:19"
$zhcn.Range("K3").Value = "2016-09-05 06:55:36"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-09-05 06:55:24"
$dede.Range("K3").Value = "2016-09-05 06:55:45"
